$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tpDictionary")

# Update the gender unit identifier from "Organism|Gender" to "Gender"
$ws.Range("F12").Value = "Gender"

# Move selection to D13 as a result of the edit
$ws.Range("D13").Select()
